# Apply the 2022-07-14 20:58:51 re-crawl update:
#  - several product rows were re-ordered by the crawler (their full
#    row contents moved to a different row), so swap/rotate the affected
#    rows back into their new positions
#  - every data row's timestamp (column O) is bumped to the new crawl time

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-07-14 20:58:51"

# Columns that hold free-text / numeric-looking text (must stay text, not
# get auto-coerced into a Double by the Excel value-assignment heuristics).
$textCols = @("A","B","C","D","G","H","I","J","K","L","M","N","O")
# Columns that are genuinely numeric (ratingAmount / ratingValue).
$numCols  = @("E","F")

function Get-RowValues($row) {
    # Read the whole row A:O back as a 1-based 2D array (Value2 avoids the
    # COM Value accessor's extra currency/date wrapping).
    return $ws.Range("A$row`:O$row").Value2
}

function Set-RowValues($row, $vals) {
    # Make sure numeric-looking text (ids, prices, ...) round-trips as text
    # by forcing a text number format on the text columns before writing.
    foreach ($col in $textCols) {
        $ws.Range("$col$row").NumberFormat = "@"
    }
    foreach ($col in $numCols) {
        $ws.Range("$col$row").NumberFormat = "General"
    }
    $ws.Range("A$row`:O$row").Value = $vals
}

# --- Row 3 / 5 / 6 rotate: 6 -> 3, 3 -> 5, 5 -> 6 ---------------------------
$row3 = Get-RowValues 3
$row5 = Get-RowValues 5
$row6 = Get-RowValues 6

Set-RowValues 3 $row6
Set-RowValues 5 $row3
Set-RowValues 6 $row5

# --- Row 13 / 14 swap --------------------------------------------------------
$row13 = Get-RowValues 13
$row14 = Get-RowValues 14

Set-RowValues 13 $row14
Set-RowValues 14 $row13

# --- Row 18 / 19 swap --------------------------------------------------------
$row18 = Get-RowValues 18
$row19 = Get-RowValues 19

Set-RowValues 18 $row19
Set-RowValues 19 $row18

# --- Row 24 / 25 swap --------------------------------------------------------
$row24 = Get-RowValues 24
$row25 = Get-RowValues 25

Set-RowValues 24 $row25
Set-RowValues 25 $row24

# --- Bump the crawl timestamp for every data row (2 through 34) ------------
for ($r = 2; $r -le 34; $r++) {
    $ws.Range("O$r").NumberFormat = "@"
    $ws.Range("O$r").Value = $newTimestamp
}
